$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = -12.091
    "B3" = 5.468000000000001
    "D3" = -7.256
    "D12" = -7.343000000000001
    "B14" = 5.63
    "B16" = 5.394
    "C18" = -12.543
    "B21" = 9.460000000000001
    "B23" = 8.000999999999999
    "C24" = -12.737
    "D24" = -8.026999999999999
    "B25" = 5.468999999999999
    "C25" = -13.597
    "D25" = -8.307
    "B26" = 5.797
    "C27" = -13.203
    "B29" = 5.372999999999999
    "C30" = -12.7
    "C31" = -13.257
    "C39" = -12.717
    "B40" = 9.277000000000001
    "D41" = -8.130000000000001
    "C42" = -12.67
    "C48" = -11.374
    "D50" = -8.045999999999999
    "C51" = -11.889
    "C52" = -11.221
    "B53" = 6.383
    "D53" = -7.053999999999999
    "C55" = -13.748
    "C56" = -12.682
    "D56" = -8.087999999999999
    "B57" = 5.178
    "C57" = -13.716
    "D57" = -8.488
    "D58" = -8.318
    "B59" = 4.811
    "C60" = -12.752
    "D61" = -7.722
    "D63" = -7.885999999999998
    "D64" = -7.495
    "B65" = 5.943
    "B69" = 5.319999999999999
    "D70" = -7.177000000000001
    "D72" = -7.373
    "C73" = -13.006
    "C74" = -12.453
    "B79" = 5.541000000000001
    "B83" = 5.475
    "D86" = -8.190999999999999
    "C89" = -10.981
    "D89" = -6.412000000000001
    "C90" = -13.298
    "B91" = 5.545
    "C92" = -11.241
    "B93" = 5.438
    "D98" = -8.612
    "B100" = 4.929999999999999
    "D100" = -8.458
    "D102" = -7.788000000000001
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
